$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "Legislature" (row 11) and shift everything below up.
$ws.Rows.Item(11).Delete()
